# Update the "想去人数" (F column) figures that were refreshed when the
# gh-pages data output was regenerated (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 2168
$ws.Range("F15").Value = 4191
$ws.Range("F18").Value = 1148
$ws.Range("F20").Value = 444
$ws.Range("F21").Value = 6295
$ws.Range("F27").Value = 1967
$ws.Range("F31").Value = 26
$ws.Range("F34").Value = 60
$ws.Range("F42").Value = 1173
$ws.Range("F45").Value = 1065
$ws.Range("F48").Value = 180

# Sheet "演出" (performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F22").Value = 180

# Sheet "本地生活" (local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value  = 480
$ws.Range("F9").Value  = 979
$ws.Range("F10").Value = 1115
$ws.Range("F11").Value = 1323
$ws.Range("F12").Value = 1620
$ws.Range("F13").Value = 29

# Sheet "全部类型" (all types - aggregated view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 480
$ws.Range("F8").Value  = 979
$ws.Range("F17").Value = 1620
$ws.Range("F18").Value = 4191
$ws.Range("F23").Value = 1148
$ws.Range("F25").Value = 444
$ws.Range("F26").Value = 6295
$ws.Range("F29").Value = 1967
$ws.Range("F32").Value = 26
$ws.Range("F35").Value = 60
$ws.Range("F43").Value = 1173
$ws.Range("F48").Value = 180
